# Vguard_Automation.xlsx - "Observation Near Miss CAPA end to end flow Test script added"
#
# Adds a new worksheet "ObservationNearMissAction" (a copy of "IncidentAction")
# right after "IncidentAction", tweaks a couple of cell values on both sheets,
# nudges the selection on "EHSObservationAction", and leaves "IncidentAction"
# as the active / selected tab when the workbook is saved.

$wb = $excel.ActiveWorkbook

# 1. Duplicate "IncidentAction" -> "ObservationNearMissAction", placed right
#    after the source sheet (same relative position the diff shows).
$incidentAction = $wb.Worksheets.Item("IncidentAction")
$incidentAction.Copy($null, $incidentAction)
$newSheet = $wb.Worksheets.Item("IncidentAction (2)")
$newSheet.Name = "ObservationNearMissAction"

# 2. Fix up the two cells that differ from the source sheet.
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("C2").Value = "Action Item of EHS"

# 3. The source "IncidentAction" sheet's header cell also changes.
$incidentAction.Range("A1").Value = "UserName"

# 4. "EHSObservationAction" gets a new selected cell (C2).
$obsAction = $wb.Worksheets.Item("EHSObservationAction")
$obsAction.Range("C2").Select()

# 5. Leave "IncidentAction" as the active sheet/tab on save.
$incidentAction.Activate()
$incidentAction.Range("A1").Select()
